$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 18: replace the old array-SUM formula with a non-array SUBSTITUTE
#     formula that errors out (#VALUE!) because SUBSTITUTE doesn't accept
#     an array range as its first argument outside of an array context.
$ws.Range("F18").Formula = '=SUBSTITUTE(F2:F17,"0","Do",1)'

# --- Row 19 (new): array-entered SUBSTITUTE formula, evaluates to text "7200"
$ws.Range("F19").FormulaArray = '=SUBSTITUTE(F3:F18,"10","Do",1)'

# --- Column widths (F, I, J) matching the authored widths as closely as
#     this host's width-rounding allows.
$ws.Columns.Item(6).ColumnWidth = 16.8333333333333
$ws.Columns.Item(9).ColumnWidth = 9.5
$ws.Columns.Item(10).ColumnWidth = 6.66666666666667

# --- Selection moves to F18:F19 with F18 active, matching the new data.
$null = $ws.Range("F18:F19").Select()
